$d = $word.ActiveDocument

function Replace-Unique($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Date heading
Replace-Unique "2025-04-27 Sunday" "2025-04-28 Monday"

# Simple unique division problems (old text occurs exactly once in the doc)
Replace-Unique "11÷7=" "53÷9="
Replace-Unique "37÷2=" "19÷6="
Replace-Unique "20÷8=" "24÷9="
Replace-Unique "17÷7=" "46÷2="
Replace-Unique "12÷5=" "24÷5="
Replace-Unique "26÷3=" "52÷4="
Replace-Unique "95÷3=" "45÷8="
Replace-Unique "23÷7=" "59÷2="
Replace-Unique "97÷3=" "55÷4="
Replace-Unique "45÷4=" "40÷4="
Replace-Unique "67÷3=" "10÷2="
Replace-Unique "81÷7=" "27÷8="
Replace-Unique "78÷6=" "66÷9="
Replace-Unique "49÷2=" "96÷9="
Replace-Unique "80÷6=" "53÷5="
Replace-Unique "94÷6=" "77÷7="
Replace-Unique "98÷2=" "97÷8="
Replace-Unique "53÷3=" "97÷2="
Replace-Unique "78÷9=" "43÷9="
Replace-Unique "68÷2=" "29÷2="
Replace-Unique "33÷4=" "90÷6="
Replace-Unique "94÷5=" "81÷9="
Replace-Unique "14÷2=" "90÷4="

# "77÷4=" appears twice in the table; handled positionally by cell.
$tbl = $d.Tables.Item(1)
$tbl.Cell(9, 3).Range.Text = "26÷3="
$tbl.Cell(17, 1).Range.Text = "91÷6="
